# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the c2e86c96 handback row on the zh-cn and de-de sheets,
# reflecting a newly generated handback report.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-23 03:09:50"
$wsZhCn.Range("H2").Value = "2016-03-23 03:10:15"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-23 03:09:54"
$wsDeDe.Range("H2").Value = "2016-03-23 03:10:22"
